$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- Update agenda rows (date + time of the three scheduled meetings) ---
# Date: 2024-02-28 -> 2024-03-11 (serial 45350 -> 45362)
$ws.Range("A2:A4").Value = 45362
# Time: 22:00:00 -> 22:00:28 (serial 0.91666666666666663 -> 0.91699074074074072)
$ws.Range("B2:B4").Value = 0.91699074074074072

# Match the "hora" header formatting (font colour + number format) to the data
# cells below it, now that it shows seconds as well.
$blackFontColor = $ws.Range("B2").Font.Color
$ws.Range("B1").Font.Color = $blackFontColor
$ws.Range("B1").NumberFormat = "h:mm:ss"

# Header row is slightly taller after the formatting refresh
$ws.Rows.Item(1).RowHeight = 19.5

# Leave the selection on B4, matching the edited workbook's UI state
$ws.Range("B4").Select()
